$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'22.405.34"
$ws.Range("E2").Value = "  -4.68%  "
$ws.Range("D3").Value = "'1.566.60"
$ws.Range("E3").Value = "  -5.04%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D6").Value = "'291.29"
$ws.Range("E6").Value = "  -2.83%  "
$ws.Range("D7").Value = "'0.3684"
$ws.Range("E7").Value = "  -2.93%  "
$ws.Range("D8").Value = "'49.57"
$ws.Range("E8").Value = "  -1.10%  "
$ws.Range("E9").Value = "  -5.55%  "
$ws.Range("E10").Value = "  -4.57%  "
$ws.Range("E11").Value = "  -6.76%  "
$ws.Range("E12").Value = "  +0.08%  "
$ws.Range("E13").Value = "  -4.50%  "
$ws.Range("D14").Value = "'6.041"
$ws.Range("E14").Value = "  -5.72%  "
$ws.Range("D15").Value = "'6.839"
$ws.Range("E15").Value = "  -7.44%  "
$ws.Range("D16").Value = "'0.00001143"
$ws.Range("E16").Value = "  -4.44%  "
$ws.Range("D17").Value = "'1.576.82"
$ws.Range("E17").Value = "  -4.74%  "
$ws.Range("D18").Value = "'89.06"
$ws.Range("E18").Value = "  -8.57%  "
$ws.Range("D19").Value = "'0.06696"
$ws.Range("E19").Value = "  -3.86%  "
$ws.Range("E20").Value = "  +0.04%  "
$ws.Range("D21").Value = "'6.248"
$ws.Range("E21").Value = "  -7.77%  "
$ws.Range("D22").Value = "'0.5296"
$ws.Range("E22").Value = "  -8.13%  "
$ws.Range("E23").Value = "  -5.29%  "
$ws.Range("E24").Value = "  -4.17%  "
$ws.Range("D25").Value = "'22.415.13"
$ws.Range("E25").Value = "  -4.74%  "
$ws.Range("D26").Value = "'2.399"
$ws.Range("E26").Value = "  -3.94%  "
$ws.Range("D27").Value = "'2.928"
$ws.Range("E27").Value = "  +0.56%  "
$ws.Range("D28").Value = "'19.83"
$ws.Range("E28").Value = "  -5.15%  "
$ws.Range("D29").Value = "'146.64"
$ws.Range("E29").Value = "  -4.55%  "
$ws.Range("D30").Value = "'4.952"
$ws.Range("E30").Value = "  -4.76%  "
$ws.Range("D31").Value = "'124.95"
$ws.Range("E31").Value = "  -5.87%  "
$ws.Range("D32").Value = "'1.751.88"
$ws.Range("E32").Value = "  -4.57%  "
$ws.Range("D33").Value = "'6.262"
$ws.Range("E33").Value = "  -9.62%  "
$ws.Range("D34").Value = "'1.976"
$ws.Range("E34").Value = "  -6.83%  "
$ws.Range("D35").Value = "'0.9821"
$ws.Range("E35").Value = "  -3.27%  "
$ws.Range("D36").Value = "'10.35"
$ws.Range("E36").Value = "  -12.90%  "
$ws.Range("D37").Value = "'0.08420"
$ws.Range("E37").Value = "  -3.70%  "
$ws.Range("E38").Value = "  -7.43%  "
$ws.Range("D39").Value = "'0.2297"
$ws.Range("E39").Value = "  -5.87%  "
$ws.Range("D40").Value = "'5.525"
$ws.Range("E40").Value = "  -7.32%  "
$ws.Range("D41").Value = "'0.06495"
$ws.Range("E41").Value = "  -4.64%  "
$ws.Range("E42").Value = "  -10.80%  "
$ws.Range("D43").Value = "'1.245"
$ws.Range("E43").Value = "  -5.57%  "
$ws.Range("D44").Value = "'0.6376"
$ws.Range("E44").Value = "  -7.66%  "
$ws.Range("D45").Value = "'14.43"
$ws.Range("E45").Value = "  -7.05%  "
$ws.Range("D46").Value = "'0.9998"
$ws.Range("E46").Value = "  -0.10%  "
$ws.Range("D47").Value = "'0.6011"
$ws.Range("E47").Value = "  -6.26%  "
$ws.Range("D48").Value = "'3.775"
$ws.Range("E48").Value = "  -3.73%  "
$ws.Range("D49").Value = "'2.105"
$ws.Range("E49").Value = "  -6.91%  "
$ws.Range("D50").Value = "'121.33"
$ws.Range("E50").Value = "  -4.91%  "
$ws.Range("D51").Value = "'0.07269"
$ws.Range("E51").Value = "  -5.97%  "
